# ECP-1152: fixes forecast creation process; adds integ test
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (ITPR001 / ITWT001 / 2020): drop Q1 (E2) from the forecast sum.
$ws.Range("N2").Formula = "=SUM(F2:H2)"

# Row 5 (ITPR001 / ITWT005 / 2020): clear the Q1 (E5) forecast entirely.
$ws.Range("E5").Clear()

# Rows 3-4 and 6-7: re-enter the "Sum Terms" formula as one gesture across
# each contiguous block so Excel records them as shared formula families.
# Row 5 keeps its own (E-column-dropped) formula, set separately below -
# mirroring row 2 (the other "2020" row) also losing Q1.
$ws.Range("N3:N4").Formula = "=SUM(E3:H3)"
$ws.Range("N5").Formula = "=SUM(F5:H5)"
$ws.Range("N6:N7").Formula = "=SUM(E6:H6)"

# Row 7 (ITPR001 / ITWT005 / 2022): Q3 (G7) forecast now has an amount.
$ws.Range("G7").Value = 1026600

# Move the active selection to E6, matching the saved view state.
$ws.Range("E6").Select()
